$wb = $excel.ActiveWorkbook

# Rename sheets (case / accent corrections from the automated export pass)
$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Drop the obsolete "Desarquivamentos Pendentes" tab entirely
[void]$wb.Worksheets("Desarquivamentos Pendentes").Delete()
